$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows effectively get their Fecha/Volumen/Precio.../Origen columns
# rotated between rows 2, 4, 10 and rows 3, 5, 12 (values "move" between rows,
# matching the target state described by the diff).

$ws.Range("D2").Value = 44999
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 2500
$ws.Range("O2").Value = 2500
$ws.Range("P2").Value = 2500
$ws.Range("R2").Value = "Región de La Araucanía"
$ws.Range("S2").Value = 2500

$ws.Range("D3").Value = 44176
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 3000
$ws.Range("O3").Value = 3000
$ws.Range("P3").Value = 3000
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 3000

$ws.Range("D4").Value = 44175
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 5000
$ws.Range("O4").Value = 5000
$ws.Range("P4").Value = 5000
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 5000

$ws.Range("D5").Value = 44323
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 3200
$ws.Range("O5").Value = 3200
$ws.Range("P5").Value = 3200
$ws.Range("S5").Value = 3200

$ws.Range("D10").Value = 44592
$ws.Range("M10").Value = 5
$ws.Range("N10").Value = 7500
$ws.Range("O10").Value = 7500
$ws.Range("P10").Value = 7500
$ws.Range("S10").Value = 7500

$ws.Range("D12").Value = 44214
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 1800
$ws.Range("O12").Value = 1800
$ws.Range("P12").Value = 1800
$ws.Range("R12").Value = "Región de La Araucanía"
$ws.Range("S12").Value = 1800
